$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.724.76"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.858.41"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").Value = "'1.033"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'323.47"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'1.030"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.4398"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").Value = "'0.3807"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "'0.07435"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "'0.8860"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").Value = "'21.61"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.876.32"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").Value = "'5.539"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "'6.748"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "'0.07178"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "'85.33"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'0.000009096"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'1.029"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "'15.53"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "27.758.28"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").Value = "'5.307"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "'11.29"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "2.093.05"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "'2.023"
$ws.Range("E25").Value = "  +6.00%  "
$ws.Range("D26").Value = "'158.13"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "'18.80"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").Value = "'5.385"
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("D29").Value = "'1.984"
$ws.Range("E29").Value = "  +3.15%  "
$ws.Range("D30").Value = "'117.94"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "'0.09016"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").Value = "'0.7816"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").Value = "'1.216"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").Value = "'3.003"
$ws.Range("E34").Value = "  +4.61%  "
$ws.Range("D35").Value = "'4.580"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("D36").Value = "'1.031"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "'1.148"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'0.01981"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'0.05310"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").Value = "'2.860"
$ws.Range("D41").Value = "'0.5213"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").Value = "'0.1686"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").Value = "'6.881"
$ws.Range("E43").Value = "  +5.21%  "
$ws.Range("D44").Value = "'8.879"
$ws.Range("E44").Value = "  +4.47%  "
$ws.Range("D45").Value = "'110.57"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").Value = "'10.71"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").Value = "'0.06602"
$ws.Range("E47").Value = "  +4.44%  "
$ws.Range("D48").Value = "'1.032"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "'1.718"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").Value = "'0.4724"
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("D51").Value = "'1.904"
$ws.Range("E51").Value = "  +0.33%  "
